$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear stale cells whose values moved to a different column ---
$ws.Range("U8").ClearContents()
$ws.Range("U9").ClearContents()

# --- Header row (row 1): party names (A1..AA1 already existed; AB1..AI1 are new) ---
$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "CC - Canary Coalition  (Coalición Canaria, CC)"
$ws.Range("C1").Value = "CiU - Convergence and Union (Convergencia i Unió, CiU)"
$ws.Range("D1").Value = "EA - Basque Solidarity (Eusko Alkartasuna, EA)"
$ws.Range("E1").Value = "ERC - Catalonian Left Republicans (Esquerra Republicana de Cataluña, ERC)"
$ws.Range("F1").Value = "HB - Unity of the People (Herri Batasuna , HB)"
$ws.Range("G1").Value = "IU - United Left (Izquierda Unida, IU)"
$ws.Range("H1").Value = "PAR - Regionalist Party Aragon (Partido Aragones Regionalista, PAR)"
$ws.Range("I1").Value = "PNV/EAJ - Basque Nationalist Party (Partido Nacionalista Vasco , PNV/EAJ)"
$ws.Range("J1").Value = "PP - Popular Party (Partido Popular, PP)"
$ws.Range("K1").Value = "PSOE - Spanish Socialist Party (Partido Socialista Obrero Español, PSOE)"
$ws.Range("L1").Value = "UV - Valencian Union (Unid Valenciana, UV)"
$ws.Range("M1").Value = "BNG - Galician Nationalist Block  (Bloque Nacionalista Gallego, BNG)"
$ws.Range("N1").Value = "CHA - Aragonese Union (Junta Aragonesista, CHA)"
$ws.Range("O1").Value = "ICV - Initiative for Catalonia Greens (niciativa per Catalunya Verds, ICV)"
$ws.Range("P1").Value = "PA - Andalusian Party ( Partido Andalucista, PA)"
$ws.Range("Q1").Value = "NB - Navarra Yes (Nafarroa Bai, NB)"
$ws.Range("R1").Value = "UPD - Progress and Democracy Union (Unión Progreso y Democracia, UPD)"
$ws.Range("S1").Value = "AMAIUR - Amiaur (Amiaur, AMAIUR)"
$ws.Range("T1").Value = "CC - Compromis Coalition (Coalició Compromís, CC)"
$ws.Range("U1").Value = "FAC - For an Alternative Cantabria  (Foro Alternativa Cantabria, FAC)"
$ws.Range("V1").Value = "C - Citizens–Party of the Citizenry (Ciudadanos – Partido de la Ciudadanía, C)"
$ws.Range("W1").Value = "CPM - Compromise-We Can-It Is Time_x0002_Coalition"
$ws.Range("X1").Value = "ECP - Together We Can"
$ws.Range("Y1").Value = "EHB - Basque Country United"
$ws.Range("Z1").Value = "P - We Can; Podemos  (Podemos, P)"
$ws.Range("AA1").Value = "PMAEU - We Can-In Tide-Anova-United Left"
$ws.Range("AB1").Value = "DL - Democracy and Freedom"
$ws.Range("AC1").Value = "JxCat - Together for Catalonia (Junts per Catalunya, JxCat)"
$ws.Range("AD1").Value = "NAplus - Navarre Addition"
$ws.Range("AE1").Value = "PAIS - `"plus Country`""
$ws.Range("AF1").Value = "PP-FORO - People's Party FORO"
$ws.Range("AG1").Value = "PRC - Cantabrian Regionalist Party"
$ws.Range("AH1").Value = "TE - Exist Teruel"
$ws.Range("AI1").Value = "V - The Voice (Vox, V)"

# --- Apply header style (bold + border + center/top align) to the newly added header cells ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AB1:AI1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column A: year labels (rows 2-10), stored as text like the rest of the column ---
$ws.Range("A2:A10").NumberFormat = "@"
$ws.Range("A2").Value = "1993"
$ws.Range("A3").Value = "1996"
$ws.Range("A4").Value = "2000"
$ws.Range("A5").Value = "2004"
$ws.Range("A6").Value = "2008"
$ws.Range("A7").Value = "2011"
$ws.Range("A8").Value = "2015"
$ws.Range("A9").Value = "2016"
$ws.Range("A10").Value = "2019"
$ws.Range("A2:A10").Style = "Normal"

# --- Data grid (numeric seat counts) ---
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 17
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 17
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 8
$ws.Range("G5").Value = 5
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 12
$ws.Range("K5").Value = 28
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = 1
$ws.Range("Q5").Value = 1
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 7
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2
$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 7
$ws.Range("K6").Value = 19
$ws.Range("M6").Value = 2
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 1
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("B8").Value = 1
$ws.Range("E8").Value = 9
$ws.Range("G8").Value = 2
$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 76
$ws.Range("K8").Value = 46
$ws.Range("V8").Value = 37
$ws.Range("W8").Value = 9
$ws.Range("X8").Value = 9
$ws.Range("Y8").Value = 2
$ws.Range("Z8").Value = 39
$ws.Range("AA8").Value = 6
$ws.Range("B9").Value = 1
$ws.Range("E9").Value = 9
$ws.Range("I9").Value = 5
$ws.Range("J9").Value = 85
$ws.Range("K9").Value = 39
$ws.Range("V9").Value = 32
$ws.Range("W9").Value = 9
$ws.Range("X9").Value = 12
$ws.Range("Y9").Value = 2
$ws.Range("Z9").Value = 39
$ws.Range("AA9").Value = 5
$ws.Range("AB9").Value = 8
$ws.Range("B10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("T10").Value = 0
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 0
$ws.Range("X10").Value = 0
$ws.Range("Y10").Value = 0
$ws.Range("Z10").Value = 0
$ws.Range("AC10").Value = 0
$ws.Range("AD10").Value = 0
$ws.Range("AE10").Value = 0
$ws.Range("AF10").Value = 0
$ws.Range("AG10").Value = 0
$ws.Range("AH10").Value = 0
$ws.Range("AI10").Value = 0
